$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 -----------------
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$nl = [char]10
$lines = @(
    "Conversión del día 💰",
    "✅ Dólar paralelo: 68",
    "",
    "Binance",
    "✅ 1000 Bs = 7.32 = 30092.21 pesos",
    "✅ 30092.21 pesos = 7.28 = 957.56 Bs",
    "",
    "Promedio competencia",
    "✅ Tasa pesos: 20",
    "✅ Tasa Bs: 20",
    "✅ % Ganancia: 20%"
)
$texto = [string]::Join($nl, $lines)
$wsHoja1.Range("A1").Value = $texto

# --- tasas: update the rate table figures -------------------------------
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 136.65
$wsTasas.Range("O10").Value = 4112.1
$wsTasas.Range("N12").Value = 4132.5
$wsTasas.Range("O12").Value = 131.5
